$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credenciales")

$ws.Range("C2").Value = "administrador"
$ws.Range("C1").Value = "disponible"
